$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 78
$ws.Cells.Item($row, 1).Value = "Globo"
$ws.Cells.Item($row, 2).Value = "RJ TV 1"
$ws.Cells.Item($row, 3).Value = "Governo"
$ws.Cells.Item($row, 4).Value = "2025-04-08T12:52"
$ws.Cells.Item($row, 5).Value = "Negativo"
$ws.Cells.Item($row, 6).Value = "Reforma administrativa suspensa em Campos. Projeto que prevê a Reforma Administrativa na Prefeitura de Campos foi suspenso. Repórter *ao vivo* em frente à Câmara Municipal. Juiz concedeu liminar suspendendo a tramitação do projeto, que foi para a Câmara 27 de Março. Previsão é de que fosse votado hoje ou amanhã. Projeto pode causar impacto de 10% a mais na folha, que é de R$ 1 bi ao ano. Pedido de suspensão partiu de uma vereadora, alegando falta de transparência. "
